$wb = $excel.ActiveWorkbook

# Add the new "Data" worksheet as the last sheet in the workbook.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Data"

# Header row (bold, matches the style used for other sheet headers).
$headers = @("event","sex","time","value")
for ($c = 0; $c -lt $headers.Count; $c++) {
  $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}
$ws.Range("A1:D1").Font.Bold = $true

# Data rows.
$rows = @(
  @("B","F",2018,338),
  @("B","F",2019,340),
  @("B","F",2020,310),
  @("B","F",2021,0),
  @("B","M",2018,346),
  @("B","M",2019,343),
  @("B","M",2020,368),
  @("B","M",2021,0),
  @("B","T",2018,684),
  @("B","T",2019,683),
  @("B","T",2020,678),
  @("B","T",2021,0),
  @("C","F",2018,1),
  @("C","F",2019,0),
  @("C","F",2020,0),
  @("C","F",2021,0),
  @("C","M",2018,-1),
  @("C","M",2019,0),
  @("C","M",2020,1),
  @("C","M",2021,0),
  @("C","T",2018,0),
  @("C","T",2019,0),
  @("C","T",2020,1),
  @("C","T",2021,0),
  @("D","F",2018,175),
  @("D","F",2019,193),
  @("D","F",2020,176),
  @("D","F",2021,0),
  @("D","M",2018,219),
  @("D","M",2019,218),
  @("D","M",2020,191),
  @("D","M",2021,0),
  @("D","T",2018,394),
  @("D","T",2019,411),
  @("D","T",2020,367),
  @("D","T",2021,0),
  @("I","F",2018,789),
  @("I","F",2019,770),
  @("I","F",2020,759),
  @("I","F",2021,0),
  @("I","M",2018,1059),
  @("I","M",2019,996),
  @("I","M",2020,938),
  @("I","M",2021,0),
  @("I","T",2018,1848),
  @("I","T",2019,1766),
  @("I","T",2020,1697),
  @("I","T",2021,0),
  @("O","F",2018,632),
  @("O","F",2019,542),
  @("O","F",2020,525),
  @("O","F",2021,0),
  @("O","M",2018,692),
  @("O","M",2019,665),
  @("O","M",2020,684),
  @("O","M",2021,0),
  @("O","T",2018,1324),
  @("O","T",2019,1207),
  @("O","T",2020,1209),
  @("O","T",2021,0),
  @("P","F",2018,24461),
  @("P","F",2019,24782),
  @("P","F",2020,25157),
  @("P","F",2021,25525),
  @("P","M",2018,26003),
  @("P","M",2019,26496),
  @("P","M",2020,26952),
  @("P","M",2021,27384),
  @("P","T",2018,50464),
  @("P","T",2019,51278),
  @("P","T",2020,52109),
  @("P","T",2021,52909),
  @("U","F",2018,24782),
  @("U","F",2019,25157),
  @("U","F",2020,25525),
  @("U","F",2021,0),
  @("U","M",2018,26496),
  @("U","M",2019,26952),
  @("U","M",2020,27384),
  @("U","M",2021,0),
  @("U","T",2018,51278),
  @("U","T",2019,52109),
  @("U","T",2020,52909),
  @("U","T",2021,0)
)

for ($r = 0; $r -lt $rows.Count; $r++) {
  $row = $rows[$r]
  for ($c = 0; $c -lt $row.Count; $c++) {
    $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
  }
}

$ws.Range("G8").Select() | Out-Null
